$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G4 (Dropdown Values for Interests1 property) to include "Writing"
$ws.Range("G4").Value = "Sports, Music, Reading, Writing"

# Update the active selection to H7 (as reflected in the saved view state)
$ws.Range("H7").Select()
